$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '66.435.99'
$ws.Range('E2').Value = '  -0.61%  '
Set-TextValue 'D3' '3.589.29'
$ws.Range('E3').Value = '  +3.18%  '
Set-TextValue 'D4' '0.998'
$ws.Range('E4').Value = '  -0.16%  '
Set-TextValue 'D5' '608.77'
$ws.Range('E5').Value = '  +0.87%  '
Set-TextValue 'D6' '145.31'
$ws.Range('E6').Value = '  +0.55%  '
Set-TextValue 'D7' '3.587.10'
$ws.Range('E7').Value = '  +3.13%  '
$ws.Range('E8').Value = '  +0.10%  '
Set-TextValue 'D9' '0.485'
$ws.Range('E9').Value = '  +1.89%  '
Set-TextValue 'D10' '0.137'
$ws.Range('E10').Value = '  -2.06%  '
Set-TextValue 'D11' '8.04'
$ws.Range('E11').Value = '  +1.45%  '
Set-TextValue 'D12' '0.413'
$ws.Range('E12').Value = '  -0.64%  '
Set-TextValue 'D13' '4.181.64'
$ws.Range('E13').Value = '  +2.77%  '
Set-TextValue 'D14' '0.0000210'
$ws.Range('E14').Value = '  -0.95%  '
Set-TextValue 'D15' '30.48'
$ws.Range('E15').Value = '  -1.58%  '
Set-TextValue 'D16' '3.579.83'
$ws.Range('E16').Value = '  +2.99%  '
Set-TextValue 'D17' '66.442.64'
$ws.Range('E17').Value = '  -0.70%  '
Set-TextValue 'D18' '11.79'
$ws.Range('E18').Value = '  +11.59%  '
$ws.Range('E19').Value = '  -0.99%  '
Set-TextValue 'D20' '6.24'
$ws.Range('E20').Value = '  -0.16%  '
Set-TextValue 'D21' '15.03'
$ws.Range('E21').Value = '  -1.51%  '
Set-TextValue 'D22' '430.84'
$ws.Range('E22').Value = '  +0.84%  '
Set-TextValue 'D23' '0.611'
$ws.Range('E23').Value = '  +2.07%  '
Set-TextValue 'D24' '78.82'
$ws.Range('E24').Value = '  -0.52%  '
Set-TextValue 'D25' '3.723.27'
$ws.Range('E25').Value = '  +3.01%  '
$ws.Range('E26').Value = '  +0.06%  '
Set-TextValue 'D27' '0.0000122'
$ws.Range('E27').Value = '  +5.48%  '
Set-TextValue 'D28' '8.13'
$ws.Range('E28').Value = '  +0.92%  '
Set-TextValue 'D29' '9.30'
$ws.Range('E29').Value = '  -3.64%  '
Set-TextValue 'D30' '2.53'
$ws.Range('E30').Value = '  +1.95%  '
Set-TextValue 'D31' '1.00'
$ws.Range('E31').Value = '  +0.02%  '
Set-TextValue 'D32' '1.50'
$ws.Range('E32').Value = '  -2.26%  '
Set-TextValue 'D33' '0.160'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D34' '25.61'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue 'D35' '3.576.21'
$ws.Range('E35').Value = '  +2.82%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '1.77'
$ws.Range('E36').Value = '  +0.85%  '
$ws.Range('B37').Value = 'USDe'
$ws.Range('C37').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D37' '1.00'
$ws.Range('E37').Value = '  -0.03%  '
Set-TextValue 'D38' '7.90'
$ws.Range('E38').Value = '  +0.26%  '
Set-TextValue 'D39' '5.69'
$ws.Range('E39').Value = '  +0.55%  '
Set-TextValue 'D40' '0.997'
$ws.Range('E40').Value = '  -0.24%  '
Set-TextValue 'D41' '171.44'
$ws.Range('E41').Value = '  -1.99%  '
Set-TextValue 'D42' '0.0863'
$ws.Range('E42').Value = '  -2.57%  '
Set-TextValue 'D43' '5.31'
$ws.Range('E43').Value = '  +0.48%  '
Set-TextValue 'D44' '0.900'
$ws.Range('E44').Value = '  +1.05%  '
Set-TextValue 'D45' '1.92'
$ws.Range('E45').Value = '  -1.68%  '
Set-TextValue 'D46' '45.94'
$ws.Range('E46').Value = '  -0.52%  '
Set-TextValue 'D47' '1.23'
$ws.Range('E47').Value = '  +2.86%  '
Set-TextValue 'D48' '26.15'
$ws.Range('E48').Value = '  -5.04%  '
Set-TextValue 'D49' '2.42'
$ws.Range('E49').Value = '  +2.45%  '
Set-TextValue 'D50' '7.18'
$ws.Range('E50').Value = '  -1.31%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '23.25'
$ws.Range('E51').Value = '  +11.79%  '
